$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "43.144.88"
$ws.Cells.Item(2, 5).Value = "  +1.02%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.384.03"
$ws.Cells.Item(3, 5).Value = "  +3.30%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.01%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "302.84"
$ws.Cells.Item(5, 5).Value = "  +0.34%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "97.00"
$ws.Cells.Item(6, 5).Value = "  +1.19%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.505"
$ws.Cells.Item(7, 5).Value = "  -0.37%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +1.12%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "34.30"
$ws.Cells.Item(10, 5).Value = "  -0.62%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.0789"
$ws.Cells.Item(11, 5).Value = "  +0.49%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  +2.32%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "18.27"
$ws.Cells.Item(13, 5).Value = "  -4.83%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "6.81"
$ws.Cells.Item(14, 5).Value = "  +1.36%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.754.17"
$ws.Cells.Item(15, 5).Value = "  +3.41%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.362.42"
$ws.Cells.Item(16, 5).Value = "  +2.32%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  +3.20%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "43.125.42"
$ws.Cells.Item(18, 5).Value = "  +1.16%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "12.19"
$ws.Cells.Item(19, 5).Value = "  -0.94%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +4.27%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.0₃0888"
$ws.Cells.Item(21, 5).Value = "  -0.21%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "68.30"
$ws.Cells.Item(22, 5).Value = "  +0.54%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "235.71"
$ws.Cells.Item(23, 5).Value = "  +0.04%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "2.24"
$ws.Cells.Item(24, 5).Value = "  -2.71%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +1.62%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.01%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "25.03"
$ws.Cells.Item(27, 5).Value = "  +3.17%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  +0.34%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "9.24"
$ws.Cells.Item(29, 5).Value = "  +1.54%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "31.33"
$ws.Cells.Item(30, 5).Value = "  -3.07%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.05%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "5.07"
$ws.Cells.Item(32, 5).Value = "  +1.14%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "0.0752"
$ws.Cells.Item(33, 5).Value = "  +7.36%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "17.40"
$ws.Cells.Item(34, 5).Value = "  -1.59%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "EnergySwap"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(35, 4).Value = "24.01"
$ws.Cells.Item(35, 5).Value = "  +18.15%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "1.87"
$ws.Cells.Item(36, 5).Value = "  +7.40%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Kaspa"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(37, 4).Value = "0.105"
$ws.Cells.Item(37, 5).Value = "  +5.49%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "WEMIXToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(38, 4).Value = "2.32"
$ws.Cells.Item(38, 5).Value = "  -0.33%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "RenderToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(39, 4).Value = "4.34"
$ws.Cells.Item(39, 5).Value = "  -2.73%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "LidoDAOToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(40, 4).Value = "2.81"
$ws.Cells.Item(40, 5).Value = "  +4.13%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -0.21%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "109.17"
$ws.Cells.Item(42, 5).Value = "  -34.18%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.952.38"
$ws.Cells.Item(43, 5).Value = "  -0.16%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "0.0281"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +2.14%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "2.75"
$ws.Cells.Item(46, 5).Value = "  +0.30%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -12.51%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "2.606.70"
$ws.Cells.Item(48, 5).Value = "  +2.94%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "Stacks"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(49, 4).Value = "1.52"
$ws.Cells.Item(49, 5).Value = "  +2.45%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "MultiversX"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(50, 4).Value = "52.46"
$ws.Cells.Item(50, 5).Value = "  -1.93%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "TrustWalletToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(51, 4).Value = "1.15"
$ws.Cells.Item(51, 5).Value = "  +1.31%  "
